$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Public ID" is an auto-incrementing identifier; clear the stored sample
# values for rows 2-5 (Theme 1-4) so they get regenerated on next insert.
$ws.Range("A2:A5").ClearContents()

# Remove sample rows for Theme 5-10 entirely (rows 6-11), leaving only
# Theme 1-4 as sample data.
$ws.Range("A6:C11").ClearContents()
$ws.Range("A6:C11").NumberFormat = "General"

# Delete the now-unused trailing empty rows 17-22 so the sheet shrinks
# back down to match the reduced sample data set.
$ws.Range("A17:H22").EntireRow.Delete()
